# "Generate Report for Handback"
#
# Two e2e files (b3493285-... and f9edfca9-...) finish handback:
#   - Overview sheet: zh-cn/de-de status columns flip from
#     "Ready for handoff" -> "Handed back: in sync with en-US"
#   - zh-cn / de-de detail sheets (rows 4 & 5):
#       Status            -> "Handed back: in sync with en-US"
#       Latest Target File-> source md file (now a hyperlink, like rows 2/3)
#       Latest Handback File -> the generated .xlf handback file
#       Latest Handback DateTime -> real timestamp (was the zero-date sentinel)

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $statusHandedBack
$overview.Range("F4").Value = $statusHandedBack
$overview.Range("E5").Value = $statusHandedBack
$overview.Range("F5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $statusHandedBack
$zhcn.Range("J4").Value = "b3493285-5a46-445e-be12-8360d30c8444.c11c6ce2dab30c71faa1cb779d4133dbbc4ddfe0.zh-cn.xlf"
$zhcn.Range("K4").Value = "2016-08-21 16:39:31"
$zhcn.Hyperlinks.Add($zhcn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e573af7a505bad92fbeb7d367d9d5255e52ef6e8/e2e/b3493285-5a46-445e-be12-8360d30c8444.md", [Type]::Missing, [Type]::Missing, "b3493285-5a46-445e-be12-8360d30c8444.md")
$zhcn.Range("I4").Font.Underline = $true
$zhcn.Range("I4").Font.Color = 15570276
$zhcn.Range("I4").Font.Name = "Calibri"
$zhcn.Range("I4").Font.Size = 11

$zhcn.Range("C5").Value = $statusHandedBack
$zhcn.Range("J5").Value = "f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.500cf196655d5d4d732f14ebffb8b9e453f2f2f5.zh-cn.xlf"
$zhcn.Range("K5").Value = "2016-08-21 16:39:31"
$zhcn.Hyperlinks.Add($zhcn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e573af7a505bad92fbeb7d367d9d5255e52ef6e8/e2e/f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.md", [Type]::Missing, [Type]::Missing, "f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.md")
$zhcn.Range("I5").Font.Underline = $true
$zhcn.Range("I5").Font.Color = 15570276
$zhcn.Range("I5").Font.Name = "Calibri"
$zhcn.Range("I5").Font.Size = 11

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $statusHandedBack
$dede.Range("J4").Value = "b3493285-5a46-445e-be12-8360d30c8444.c11c6ce2dab30c71faa1cb779d4133dbbc4ddfe0.de-de.xlf"
$dede.Range("K4").Value = "2016-08-21 16:39:37"
$dede.Hyperlinks.Add($dede.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5ba7c0e068b00336208b27de7812b5d74602e418/e2e/b3493285-5a46-445e-be12-8360d30c8444.md", [Type]::Missing, [Type]::Missing, "b3493285-5a46-445e-be12-8360d30c8444.md")
$dede.Range("I4").Font.Underline = $true
$dede.Range("I4").Font.Color = 15570276
$dede.Range("I4").Font.Name = "Calibri"
$dede.Range("I4").Font.Size = 11

$dede.Range("C5").Value = $statusHandedBack
$dede.Range("J5").Value = "f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.500cf196655d5d4d732f14ebffb8b9e453f2f2f5.de-de.xlf"
$dede.Range("K5").Value = "2016-08-21 16:39:37"
$dede.Hyperlinks.Add($dede.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5ba7c0e068b00336208b27de7812b5d74602e418/e2e/f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.md", [Type]::Missing, [Type]::Missing, "f9edfca9-e6b5-4497-a6e5-27bb7f7ec477.md")
$dede.Range("I5").Font.Underline = $true
$dede.Range("I5").Font.Color = 15570276
$dede.Range("I5").Font.Name = "Calibri"
$dede.Range("I5").Font.Size = 11

Write-Output "Handback report generated"
